# Release-Notes.xlsx update:
# A new/updated "Implement CI-CD with GitHub Actions" entry now has the most
# recent timestamp among the folder list (after the newest "Azure_Well-
# Architected..." row), so it moves to row 3 of the "Folder Inventory" sheet,
# pushing every row from the old row 3 through the old row 62 down by one
# (the entry's old position was row 63 before the update). Rows 64+ are
# untouched. The Metadata sheet's "Generated On" and "Workflow Run" values
# are also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folder Inventory")

# Shift rows 62 .. 3 down to 63 .. 4 (work bottom-up so we don't clobber
# source data before it is read). NOTE: use Value2 for reads -- the plain
# Value getter in this host mis-resolves to the property descriptor.
for ($r = 62; $r -ge 3; $r--) {
    $dest = $r + 1
    $name = $ws.Cells.Item($r, 1).Value2
    $updated = $ws.Cells.Item($r, 3).Value2
    $count = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($dest, 1).Value2 = $name
    $ws.Cells.Item($dest, 2).Value2 = $name
    $ws.Cells.Item($dest, 3).Value2 = $updated
    $ws.Cells.Item($dest, 4).Value2 = $count
    $ws.Cells.Item($dest, 5).Value2 = "Root"
}

# Row 3 becomes the refreshed "Implement CI-CD with GitHub Actions" entry.
$ws.Cells.Item(3, 1).Value2 = "Implement CI-CD with GitHub Actions"
$ws.Cells.Item(3, 2).Value2 = "Implement CI-CD with GitHub Actions"
$ws.Cells.Item(3, 3).Value2 = "2025-06-13 15:19:07 +0000"
$ws.Cells.Item(3, 4).Value2 = 1
$ws.Cells.Item(3, 5).Value2 = "Root"

# Metadata sheet refresh.
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value2 = "2025-06-13 15:19:23 UTC"
# "Workflow Run" is stored as text (not a number) in the workbook, so force
# text storage with a leading quote-prefix -- otherwise the numeric-looking
# "6" gets auto-coerced into a real number by Excel's type inference.
$meta.Cells.Item(5, 2).Value2 = "'6"
